$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily reading (2026/01/30, 金, time=2, ranking=201) was recorded and
# inserted right after the existing 2026/01/29 entries (row 746), pushing the
# 2026/12/29-2027/01/05 block (old rows 747-788) down by one row.
$ws.Rows(747).Insert()

# Fill in the newly inserted row. Force column A to be stored as literal
# text (matching every other date cell in the sheet) instead of letting
# Excel auto-convert the "yyyy/mm/dd" string into a date serial number.
$ws.Range("A747").NumberFormat = "@"
$ws.Range("A747").Value = "2026/01/30"
$ws.Range("A747").Style = "Normal"
$ws.Range("B747").Value = "金"
$ws.Range("C747").Value = 2
$ws.Range("D747").Value = 201
